$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "grouping1" (AA) and "grouping2" (AB) columns entirely,
# shifting all following columns (tissue, disease, etc.) left.
$ws.Range("AA:AB").Delete()
